$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell "Save" in H1, copying the style/format from G1 (the "sum" header)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the data cell H2 with value 0 (no special style, matching other data cells)
$ws.Range("H2").Value = 0
